{"js": "// Update the two-digit division worksheet: replace each \"A\u00f7B=\" prompt\n// in the table with its new value. Old values are unique in the\n// document, so a simple search/replace per pair is safe. The pairs are\n// applied in document order so that a value which is simultaneously an\n// old prompt (earlier in the doc) and a new prompt (later in the doc) -\n// here \"85\u00f72=\" - is not re-matched after it is written.\nconst replacements = [\n  [\"71\u00f78=\", \"19\u00f76=\"],\n  [\"91\u00f77=\", \"34\u00f79=\"],\n  [\"56\u00f78=\", \"20\u00f78=\"],\n  [\"52\u00f76=\", \"74\u00f73=\"],\n  [\"86\u00f72=\", \"60\u00f73=\"],\n  [\"88\u00f74=\", \"15\u00f72=\"],\n  [\"11\u00f76=\", \"12\u00f79=\"],\n  [\"86\u00f73=\", \"29\u00f74=\"],\n  [\"21\u00f79=\", \"92\u00f76=\"],\n  [\"88\u00f79=\", \"69\u00f73=\"],\n  [\"71\u00f73=\", \"78\u00f73=\"],\n  [\"93\u00f74=\", \"58\u00f77=\"],\n  [\"33\u00f78=\", \"93\u00f75=\"],\n  [\"49\u00f77=\", \"27\u00f78=\"],\n  [\"90\u00f75=\", \"50\u00f72=\"],\n  [\"45\u00f72=\", \"23\u00f78=\"],\n  [\"85\u00f72=\", \"40\u00f76=\"],\n  [\"61\u00f78=\", \"47\u00f78=\"],\n  [\"70\u00f79=\", \"27\u00f77=\"],\n  [\"94\u00f75=\", \"53\u00f75=\"],\n  [\"92\u00f72=\", \"69\u00f76=\"],\n  [\"65\u00f72=\", \"15\u00f78=\"],\n  [\"78\u00f74=\", \"96\u00f74=\"],\n  [\"23\u00f77=\", \"85\u00f72=\"],\n  [\"98\u00f78=\", \"84\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit division worksheet: replace each \"A\u00f7B=\" prompt\n# in the table with its new value. Old values are unique in the\n# document, so Find/Replace per pair is safe. The pairs are applied in\n# document order so that a value which is simultaneously an old prompt\n# (earlier in the doc) and a new prompt (later in the doc) - here\n# \"85\u00f72=\" - is not re-matched after it has already been written.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"71\u00f78=\", \"19\u00f76=\"),\n    @(\"91\u00f77=\", \"34\u00f79=\"),\n    @(\"56\u00f78=\", \"20\u00f78=\"),\n    @(\"52\u00f76=\", \"74\u00f73=\"),\n    @(\"86\u00f72=\", \"60\u00f73=\"),\n    @(\"88\u00f74=\", \"15\u00f72=\"),\n    @(\"11\u00f76=\", \"12\u00f79=\"),\n    @(\"86\u00f73=\", \"29\u00f74=\"),\n    @(\"21\u00f79=\", \"92\u00f76=\"),\n    @(\"88\u00f79=\", \"69\u00f73=\"),\n    @(\"71\u00f73=\", \"78\u00f73=\"),\n    @(\"93\u00f74=\", \"58\u00f77=\"),\n    @(\"33\u00f78=\", \"93\u00f75=\"),\n    @(\"49\u00f77=\", \"27\u00f78=\"),\n    @(\"90\u00f75=\", \"50\u00f72=\"),\n    @(\"45\u00f72=\", \"23\u00f78=\"),\n    @(\"85\u00f72=\", \"40\u00f76=\"),\n    @(\"61\u00f78=\", \"47\u00f78=\"),\n    @(\"70\u00f79=\", \"27\u00f77=\"),\n    @(\"94\u00f75=\", \"53\u00f75=\"),\n    @(\"92\u00f72=\", \"69\u00f76=\"),\n    @(\"65\u00f72=\", \"15\u00f78=\"),\n    @(\"78\u00f74=\", \"96\u00f74=\"),\n    @(\"23\u00f77=\", \"85\u00f72=\"),\n    @(\"98\u00f78=\", \"84\u00f78=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]$wdFindContinue, [ref]$false, [ref]$find.Replacement.Text, [ref]$wdReplaceAll)\n}\n"}
